{"js": "// 1) \"\uc9c0\ub3c4\uc5d0 \ub0b4 \uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\" -> \"\uc9c0\ub3c4\uc5d0 \ub0b4\uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\"\n{\n  const results = context.document.body.search(\"\uc9c0\ub3c4\uc5d0 \ub0b4 \uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\", {\n    matchCase: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text for edit #1 (\ub0b4 \uc704\uce58 \ud45c\uc2dc)\");\n  }\n  results.items[0].insertText(\"\uc9c0\ub3c4\uc5d0 \ub0b4\uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c\" ->\n//    \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c \ucc3e\uae30\"\n{\n  const results = context.document.body.search(\n    \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text for edit #2 (\uc9d0 \ubcf4\uad00\uc18c)\");\n  }\n  results.items[0].insertText(\n    \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c \ucc3e\uae30\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3) \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9\" ->\n//    \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9 ~ \ucd9c\ubc1c:\ub0b4\uc704\uce58&\ub3c4\ucc29:\uc2dc\uc7a5\uc8fc\uc18c, \uac00\ub294\ubc29\ubc95, \uad50\ud1b5\uc218\ub2e8\ubcc4 \uc18c\uc694\uc2dc\uac04\"\n{\n  const results = context.document.body.search(\"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text for edit #3 (\ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9)\");\n  }\n  results.items[0].insertText(\n    \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9 ~ \ucd9c\ubc1c:\ub0b4\uc704\uce58&\ub3c4\ucc29:\uc2dc\uc7a5\uc8fc\uc18c, \uac00\ub294\ubc29\ubc95, \uad50\ud1b5\uc218\ub2e8\ubcc4 \uc18c\uc694\uc2dc\uac04\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 4) Insert a new block of 20 paragraphs right after the \"2-3 \uba39\uac70\ub9ac&\ubcfc\uac70\ub9ac\" paragraph\n//    (and before the following paragraph that holds a single space \" \").\n{\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  let anchor = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === \"2-3 \uba39\uac70\ub9ac&\ubcfc\uac70\ub9ac\") {\n      anchor = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!anchor) {\n    throw new Error(\"Could not find anchor paragraph '2-3 \uba39\uac70\ub9ac&\ubcfc\uac70\ub9ac'\");\n  }\n\n  // null => insert a paragraph with no text run (matches the blank separator\n  // paragraphs that appear in the source diff).\n  const newParagraphTexts = [\n    \"(1)\uba39\uac70\ub9ac\",\n    \"\uc778\uae30\uba39\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc704\uce58\uae30\ubc18 \ucd94\ucc9c \uba39\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc0ac\uc6a9\uc790 \uc815\ubcf4\uae30\ubc18 \ucd94\ucc9c \uba39\uac70\ub9ac\",\n    \"(2)\ubcfc\uac70\ub9ac\",\n    \"\uc778\uae30\ubcfc\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc704\uce58\uae30\ubc18 \ucd94\ucc9c \ubcfc\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc0ac\uc6a9\uc790 \uc815\ubcf4\uae30\ubc18 \ucd94\ucc9c \ubcfc\uac70\ub9ac\",\n    \"2-4 \uc8fc\ucc28\",\n    \"(1)\uacf5\uc601\uc8fc\ucc28\uc7a5\",\n    \"\uc9c0\ub3c4\uc5d0 \ub0b4\uc704\uce58 \uc548\ub0b4&\uc8fc\uc704 \uacf5\uc601\uc8fc\uc790\ucc3d \uc548\ub0b4\",\n    \"(2)\uc2dc\uc7a5\uc8fc\ucc28\uc7a5\",\n    \"\uc9c0\ub3c4\uc5d0 \uc8fc\ucc28\uc7a5\uc704\uce58 \uc548\ub0b4&\uc8fc\ucc28\uc694\uae08 \uc548\ub0b4\",\n    \"2-5 \ucfe0\ud3f0\",\n    \"(1)\ub0b4 \ucfe0\ud3f0\",\n    \"\uc2dc\uc7a5\ub0b4 \uc0ac\uc6a9\uac00\ub2a5 \ucfe0\ud3f0 \ub9ac\uc2a4\ud2b8&\ub0b4 \ubcf4\uc720 \ucfe0\ud3f0 \ub9ac\uc2a4\ud2b8\",\n    \"(2)\uac8c\uc784\",\n    \"\ud37c\uc990 \ud2c0 - \ud37c\uc990\ucc3e\uae30\ub97c \ud1b5\ud574 \ud37c\uc990\uc744 \ud68d\ub4dd, \ubaa8\ub450 \ud68d\ub4dd\ud558\uba74 \ud574\ub2f9 \uc2dc\uc7a5\uc774\ub098 \uc9c0\uc5ed\uc758 \ub9c8\uc2a4\ucf54\ud2b8\ub098 \uc0c1\uc9d5\uc774 \ub4f1\uc7a5(\uc2dc\uc7a5&\uc9c0\uc5ed \ud64d\ubcf4), \ucfe0\ud3f0 \ud639\uc740 \ud3ec\uc778\ud2b8 \ud68d\ub4dd \",\n    null,\n    \"\ud558\ub2e8 - \ud37c\uc990\ucc3e\uae30 \ubc84\ud2bc(\uce74\uba54\ub77c), \uce74\uba54\ub77c\ub85c \uc9c0\uc815\ub41c \uc704\uce58\uc5d0 \uc788\ub294 \ud37c\uc990 \uc870\uac01\uc744 \ucc3e\uae30, \ud37c\uc990\uc744 \ud074\ub9ad\ud558\uba74 \ud37c\uc990 \uc870\uac01\uc774 \ub9de\ucdb0\uc9d0\",\n    null,\n    \"2-6 \uc0c1\uc138\uc815\ubcf4\",\n    \"\uc2dc\uc7a5 \uc6d4\ubcc4 \uc8fc\uc694\uc77c\uc815 \uc548\ub0b4(\uc601\uc5c5\uc2dc\uac04, \ud734\ubb34\uc77c, \ud589\uc0ac)\",\n    null,\n  ];\n\n  let current = anchor;\n  for (const text of newParagraphTexts) {\n    current = current.insertParagraph(text === null ? \"\" : text, Word.InsertLocation.after);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) \"\uc9c0\ub3c4\uc5d0 \ub0b4 \uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\" -> \"\uc9c0\ub3c4\uc5d0 \ub0b4\uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$ok1 = $find.Execute(\n    \"\uc9c0\ub3c4\uc5d0 \ub0b4 \uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\", $false, $false, $false, $false, $false, $true,\n    0, $false, \"\uc9c0\ub3c4\uc5d0 \ub0b4\uc704\uce58 \ud45c\uc2dc & \uc704\uce58 \uc8fc\ubcc0 \uc2dc\uc7a5 \uc548\ub0b4\", 2\n)\nif (-not $ok1) { throw \"Edit #1 failed: target text not found (\ub0b4 \uc704\uce58 \ud45c\uc2dc)\" }\n\n# --- 2) \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c\" ->\n#        \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c \ucc3e\uae30\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$ok2 = $find2.Execute(\n    \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c\", $false, $false, $false, $false, $false, $true,\n    0, $false, \"- \uc2dc\uc7a5 \ub0b4 \uc9c0\ub3c4, \uc704\uce58 \uae30\ub85d&\uc704\uce58 \uae30\ub85d \ubaa9\ub85d, \ub85c\ub4dc\ub9f5(VR), \ud654\uc7a5\uc2e4, \uc9d0 \ubcf4\uad00\uc18c \ucc3e\uae30\", 2\n)\nif (-not $ok2) { throw \"Edit #2 failed: target text not found (\uc9d0 \ubcf4\uad00\uc18c)\" }\n\n# --- 3) \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9\" -> \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9 ~ \ucd9c\ubc1c:\ub0b4\uc704\uce58&\ub3c4\ucc29:\uc2dc\uc7a5\uc8fc\uc18c, \uac00\ub294\ubc29\ubc95, \uad50\ud1b5\uc218\ub2e8\ubcc4 \uc18c\uc694\uc2dc\uac04\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$ok3 = $find3.Execute(\n    \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9\", $false, $false, $false, $false, $false, $true,\n    0, $false, \"- \ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9 ~ \ucd9c\ubc1c:\ub0b4\uc704\uce58&\ub3c4\ucc29:\uc2dc\uc7a5\uc8fc\uc18c, \uac00\ub294\ubc29\ubc95, \uad50\ud1b5\uc218\ub2e8\ubcc4 \uc18c\uc694\uc2dc\uac04\", 2\n)\nif (-not $ok3) { throw \"Edit #3 failed: target text not found (\ub124\uc774\ubc84\ub098 \uad6c\uae00\uacfc \uc5f0\ub3d9)\" }\n\n# --- 4) Insert a block of 20 new paragraphs right after \"2-3 \uba39\uac70\ub9ac&\ubcfc\uac70\ub9ac\"\n#        (and before the following paragraph, which holds a single space \" \").\n$rng = $d.Content\n$foundAnchor = $rng.Find.Execute(\"2-3 \uba39\uac70\ub9ac&\ubcfc\uac70\ub9ac\")\nif (-not $foundAnchor) { throw \"Edit #4 failed: anchor paragraph '2-3 \uba39\uac70\ub9ac&\ubcfc\uac70\ub9ac' not found\" }\n$rng.Collapse(0)\n\n$newLines = @(\n    \"(1)\uba39\uac70\ub9ac\",\n    \"\uc778\uae30\uba39\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc704\uce58\uae30\ubc18 \ucd94\ucc9c \uba39\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc0ac\uc6a9\uc790 \uc815\ubcf4\uae30\ubc18 \ucd94\ucc9c \uba39\uac70\ub9ac\",\n    \"(2)\ubcfc\uac70\ub9ac\",\n    \"\uc778\uae30\ubcfc\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc704\uce58\uae30\ubc18 \ucd94\ucc9c \ubcfc\uac70\ub9ac \ub9ac\uc2a4\ud2b8+\ub354\ubcf4\uae30, \uc0ac\uc6a9\uc790 \uc815\ubcf4\uae30\ubc18 \ucd94\ucc9c \ubcfc\uac70\ub9ac\",\n    \"2-4 \uc8fc\ucc28\",\n    \"(1)\uacf5\uc601\uc8fc\ucc28\uc7a5\",\n    \"\uc9c0\ub3c4\uc5d0 \ub0b4\uc704\uce58 \uc548\ub0b4&\uc8fc\uc704 \uacf5\uc601\uc8fc\uc790\ucc3d \uc548\ub0b4\",\n    \"(2)\uc2dc\uc7a5\uc8fc\ucc28\uc7a5\",\n    \"\uc9c0\ub3c4\uc5d0 \uc8fc\ucc28\uc7a5\uc704\uce58 \uc548\ub0b4&\uc8fc\ucc28\uc694\uae08 \uc548\ub0b4\",\n    \"2-5 \ucfe0\ud3f0\",\n    \"(1)\ub0b4 \ucfe0\ud3f0\",\n    \"\uc2dc\uc7a5\ub0b4 \uc0ac\uc6a9\uac00\ub2a5 \ucfe0\ud3f0 \ub9ac\uc2a4\ud2b8&\ub0b4 \ubcf4\uc720 \ucfe0\ud3f0 \ub9ac\uc2a4\ud2b8\",\n    \"(2)\uac8c\uc784\",\n    \"\ud37c\uc990 \ud2c0 - \ud37c\uc990\ucc3e\uae30\ub97c \ud1b5\ud574 \ud37c\uc990\uc744 \ud68d\ub4dd, \ubaa8\ub450 \ud68d\ub4dd\ud558\uba74 \ud574\ub2f9 \uc2dc\uc7a5\uc774\ub098 \uc9c0\uc5ed\uc758 \ub9c8\uc2a4\ucf54\ud2b8\ub098 \uc0c1\uc9d5\uc774 \ub4f1\uc7a5(\uc2dc\uc7a5&\uc9c0\uc5ed \ud64d\ubcf4), \ucfe0\ud3f0 \ud639\uc740 \ud3ec\uc778\ud2b8 \ud68d\ub4dd \",\n    \"\",\n    \"\ud558\ub2e8 - \ud37c\uc990\ucc3e\uae30 \ubc84\ud2bc(\uce74\uba54\ub77c), \uce74\uba54\ub77c\ub85c \uc9c0\uc815\ub41c \uc704\uce58\uc5d0 \uc788\ub294 \ud37c\uc990 \uc870\uac01\uc744 \ucc3e\uae30, \ud37c\uc990\uc744 \ud074\ub9ad\ud558\uba74 \ud37c\uc990 \uc870\uac01\uc774 \ub9de\ucdb0\uc9d0\",\n    \"\",\n    \"2-6 \uc0c1\uc138\uc815\ubcf4\",\n    \"\uc2dc\uc7a5 \uc6d4\ubcc4 \uc8fc\uc694\uc77c\uc815 \uc548\ub0b4(\uc601\uc5c5\uc2dc\uac04, \ud734\ubb34\uc77c, \ud589\uc0ac)\",\n    \"\"\n)\n$block = \"`r\" + ($newLines -join \"`r\")\n$rng.InsertAfter($block)\n"}
